# Apply the Gantt-Chart dates/duration edits described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 5 (Task 2): push the end/actual-finish date forward ---
$ws.Range("E5").Value = 43094
$ws.Range("F5").Value = 43094

# --- Row 6 (Task 3): start now depends on the new E5, duration/variance recompute ---
$ws.Range("D6").Formula = "=E5+1"
$ws.Range("C6").Formula = "=E6-D6"
$ws.Range("E6").Value = 43132
$ws.Range("F6").Value = 43183

# --- Row 7 (Task 4): start now depends on the new F6 ---
$ws.Range("D7").Formula = "=F6+1"
$ws.Range("C7").Formula = "=E7-D7"
$ws.Range("E7").Value = 43191
$ws.Range("F7").Value = 43191

# --- Row 8 (Task 5): start typed over as a literal date, plus new end dates ---
$ws.Range("D8").Value = 43187
$ws.Range("E8").Value = 43197
$ws.Range("F8").Value = 43197

# --- View state: scrolled down, new active selection ---
$ws.Application.ActiveWindow.ScrollRow = 65
$ws.Range("F9").Select()
